$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value
$ws.Range("B2").Value = -0.344

# Delete row 4 (shifts rows 5,6,7 up to become 4,5,6)
$ws.Rows(4).Delete()
